$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "-"
$ws.Range("B8").Value = "-"
$ws.Range("F9").Value = "-"
$ws.Range("D10").Value = "-"
